# Refresh the crypto price/volume table to match the latest GitHub Actions scrape.
# (commit: "Updated cryptos list on Fri May  5 09:00:47 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.073.60"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "1.897.66"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'325.72"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").Value = "'0.3889"
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("D9").Value = "'0.07854"
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("D10").Value = "'0.9899"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").Value = "'21.84"
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("D12").Value = "1.865.15"
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("D13").Value = "'5.771"
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "'0.07006"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "'87.85"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "'0.000009933"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").Value = "'17.01"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").Value = "'0.9993"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").Value = "29.072.94"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("D22").Value = "'5.313"
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("D23").Value = "'11.08"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "2.097.78"
$ws.Range("E24").Value = "  -1.90%  "
$ws.Range("D25").Value = "'2.098"
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("D26").Value = "'155.95"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").Value = "'19.37"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").Value = "'5.880"
$ws.Range("E28").Value = "  -3.71%  "
$ws.Range("D29").Value = "'118.55"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").Value = "'1.874"
$ws.Range("E30").Value = "  -6.02%  "
$ws.Range("D31").Value = "'0.09325"
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("D32").Value = "'0.8969"
$ws.Range("E32").Value = "  -2.65%  "
$ws.Range("D33").Value = "'5.233"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("D34").Value = "'1.320"
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("E35").Value = "  -3.58%  "
$ws.Range("D36").Value = "'0.05790"
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("D37").Value = "'1.172"
$ws.Range("E37").Value = "  -2.63%  "
$ws.Range("D38").Value = "'0.02082"
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("D39").Value = "'0.9995"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "'7.673"
$ws.Range("E40").Value = "  -3.34%  "
$ws.Range("D41").Value = "'0.5671"
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("D42").Value = "'0.1794"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").Value = "'9.706"
$ws.Range("E43").Value = "  -2.46%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'2.231"
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'11.88"
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("D46").Value = "'0.5341"
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").Value = "'0.07005"
$ws.Range("D48").Value = "'1.847"
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("D49").Value = "'2.553"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").Value = "'112.70"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").Value = "'1.046"
$ws.Range("E51").Value = "  -2.34%  "
